$wb = $excel.ActiveWorkbook

# "展览" sheet (sheet1): update "想去人数" (want-to-go count) values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 12993
$wsExhibit.Range("F10").Value = 12975
$wsExhibit.Range("F14").Value = 7728

# "全部类型" sheet (sheet4): update corresponding rows (aggregated view)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 12993
$wsAll.Range("F11").Value = 12976
$wsAll.Range("F15").Value = 7728
